$d = $word.ActiveDocument

$replacements = @(
    @("314÷7=", "369÷7="),
    @("454÷9=", "633÷5="),
    @("509÷9=", "406÷6="),
    @("319÷5=", "423÷6="),
    @("919÷3=", "177÷8="),
    @("563÷7=", "674÷2="),
    @("821÷9=", "553÷2="),
    @("271÷4=", "373÷7="),
    @("320÷6=", "312÷6="),
    @("360÷3=", "903÷3="),
    @("456÷2=", "404÷7="),
    @("816÷2=", "641÷9="),
    @("716÷2=", "433÷5="),
    @("709÷6=", "786÷9="),
    @("927÷4=", "588÷7="),
    @("157÷8=", "519÷9="),
    @("770÷6=", "157÷2="),
    @("696÷2=", "320÷6="),
    @("419÷7=", "706÷8="),
    @("474÷6=", "716÷6="),
    @("453÷4=", "136÷6="),
    @("683÷7=", "577÷4="),
    @("497÷4=", "455÷4="),
    @("385÷8=", "302÷2="),
    @("664÷8=", "948÷8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
